$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new Sprint-2 hours-log entry for row 7
$ws.Range("A7").Value = 42840
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "SF-2"
$ws.Range("D7").Value = "Added a checkWin method and refactored some of the existing code for stopping when the game is won"

# Match the formatting already used by the other logged rows (centered,
# bordered, Calibri 12) -- applying it fresh creates a new style record,
# same as the rest of this sheet's rows.
$rng = $ws.Range("C7:D7")
$rng.Font.Name = "Calibri"
$rng.Font.Size = 12
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$rng.Borders.LineStyle = 1

# Move the active selection down to A8, as happens after entering data in A7
$ws.Range("A8").Select() | Out-Null
